$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp title in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 3 de Julio de 2020 a las 12:07"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Cells.Item(4,2).Value = 2837237
$ws.Cells.Item(4,3).Value = 48
$ws.Cells.Item(4,4).Value = 1191306
$ws.Cells.Item(4,5).Value = 1514446

# Row 21: 'Banglades' -> 'Banglades'
$ws.Cells.Item(21,2).Value = 156391
$ws.Cells.Item(21,3).Value = 3114
$ws.Cells.Item(21,4).Value = 68048
$ws.Cells.Item(21,5).Value = 86375
$ws.Cells.Item(21,7).Value = 42
$ws.Cells.Item(21,8).Value = 1968

# Row 31: 'Ecuador' -> 'Indonesia'
$ws.Cells.Item(31,1).Value = "Indonesia"
$ws.Cells.Item(31,2).Value = 60695
$ws.Cells.Item(31,3).Value = 1301
$ws.Cells.Item(31,4).Value = 27568
$ws.Cells.Item(31,5).Value = 30091
$ws.Cells.Item(31,7).Value = 49
$ws.Cells.Item(31,8).Value = 3036

# Row 32: 'Indonesia' -> 'Ecuador'
$ws.Cells.Item(32,1).Value = "Ecuador"
$ws.Cells.Item(32,2).Value = 59468
$ws.Cells.Item(32,4).Value = 28032
$ws.Cells.Item(32,5).Value = 26797
$ws.Cells.Item(32,8).Value = 4639

# Row 40: 'Portugal' -> 'Oman'
$ws.Cells.Item(40,1).Value = "Oman"
$ws.Cells.Item(40,2).Value = 43929
$ws.Cells.Item(40,3).Value = 1374
$ws.Cells.Item(40,4).Value = 26169
$ws.Cells.Item(40,5).Value = 17567
$ws.Cells.Item(40,7).Value = 5
$ws.Cells.Item(40,8).Value = 193

# Row 41: 'Oman' -> 'Portugal'
$ws.Cells.Item(41,1).Value = "Portugal"
$ws.Cells.Item(41,2).Value = 42782
$ws.Cells.Item(41,4).Value = 28097
$ws.Cells.Item(41,5).Value = 13098
$ws.Cells.Item(41,8).Value = 1587

# Row 42: 'Filipinas' -> 'Filipinas'
$ws.Cells.Item(42,2).Value = 40336
$ws.Cells.Item(42,3).Value = 1531
$ws.Cells.Item(42,4).Value = 11073
$ws.Cells.Item(42,5).Value = 27983
$ws.Cells.Item(42,7).Value = 6
$ws.Cells.Item(42,8).Value = 1280

# Row 44: 'Panama' -> 'Polonia'
$ws.Cells.Item(44,1).Value = "Polonia"
$ws.Cells.Item(44,2).Value = 35405
$ws.Cells.Item(44,3).Value = 259
$ws.Cells.Item(44,4).Value = 22651
$ws.Cells.Item(44,5).Value = 11247
$ws.Cells.Item(44,7).Value = 15
$ws.Cells.Item(44,8).Value = 1507

# Row 45: 'Polonia' -> 'Panama'
$ws.Cells.Item(45,1).Value = "Panama"
$ws.Cells.Item(45,2).Value = 35237
$ws.Cells.Item(45,4).Value = 16445
$ws.Cells.Item(45,5).Value = 18125
$ws.Cells.Item(45,8).Value = 667

# Row 60: 'Austria' -> 'Austria'
$ws.Cells.Item(60,2).Value = 18050
$ws.Cells.Item(60,3).Value = 109
$ws.Cells.Item(60,4).Value = 16558
$ws.Cells.Item(60,5).Value = 787

# Row 74: 'Malasia' -> 'Malasia'
$ws.Cells.Item(74,2).Value = 8648
$ws.Cells.Item(74,3).Value = 5
$ws.Cells.Item(74,4).Value = 8446
$ws.Cells.Item(74,5).Value = 81

# Row 75: 'Australia' -> 'Australia'
$ws.Cells.Item(75,2).Value = 8255
$ws.Cells.Item(75,3).Value = 254
$ws.Cells.Item(75,4).Value = 7319
$ws.Cells.Item(75,5).Value = 832

# Row 76: 'Finlandia' -> 'Consejo Danes para los Refugiados'
$ws.Cells.Item(76,1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(76,2).Value = 7311
$ws.Cells.Item(76,3).Value = 122
$ws.Cells.Item(76,4).Value = 2684
$ws.Cells.Item(76,5).Value = 4448
$ws.Cells.Item(76,7).Value = 3
$ws.Cells.Item(76,8).Value = 179

# Row 77: 'Consejo Danes para los Refugiados' -> 'El Salvador'
$ws.Cells.Item(77,1).Value = "El Salvador"
$ws.Cells.Item(77,2).Value = 7267
$ws.Cells.Item(77,3).Value = 267
$ws.Cells.Item(77,4).Value = 4268
$ws.Cells.Item(77,5).Value = 2797
$ws.Cells.Item(77,7).Value = 11
$ws.Cells.Item(77,8).Value = 202

# Row 78: 'Senegal' -> 'Finlandia'
$ws.Cells.Item(78,1).Value = "Finlandia"
$ws.Cells.Item(78,2).Value = 7241
$ws.Cells.Item(78,4).Value = 6700
$ws.Cells.Item(78,5).Value = 213
$ws.Cells.Item(78,8).Value = 328

# Row 79: 'El Salvador' -> 'Senegal'
$ws.Cells.Item(79,1).Value = "Senegal"
$ws.Cells.Item(79,2).Value = 7054
$ws.Cells.Item(79,4).Value = 4599
$ws.Cells.Item(79,5).Value = 2334
$ws.Cells.Item(79,7).Value = 0
$ws.Cells.Item(79,8).Value = 121

# Row 140: 'Uganda' -> 'Uganda'
$ws.Cells.Item(140,2).Value = 911
$ws.Cells.Item(140,3).Value = 9
$ws.Cells.Item(140,4).Value = 849
$ws.Cells.Item(140,5).Value = 62

# Row 202: 'Nueva Caledonia' -> 'Santa Lucia'
$ws.Cells.Item(202,1).Value = "Santa Lucia"
$ws.Cells.Item(202,2).Value = 22
$ws.Cells.Item(202,3).Value = 3
$ws.Cells.Item(202,4).Value = 19
$ws.Cells.Item(202,5).Value = 3

# Row 203: 'Laos' -> 'Nueva Caledonia'
$ws.Cells.Item(203,1).Value = "Nueva Caledonia"
$ws.Cells.Item(203,2).Value = 21
$ws.Cells.Item(203,4).Value = 21

# Row 204: 'Santa Lucia' -> 'Laos'
$ws.Cells.Item(204,1).Value = "Laos"

# Row 209: 'Groenlandia' -> 'Islas Malvinas'
$ws.Cells.Item(209,1).Value = "Islas Malvinas"

# Row 210: 'Islas Malvinas' -> 'Groenlandia'
$ws.Cells.Item(210,1).Value = "Groenlandia"
